$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cell updates in the three "SmartRules" test tables ---
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 2
$ws.Range("C30").Value = 3

$ws.Range("C35").Value = 4
$ws.Range("C36").Value = 5
$ws.Range("C37").Value = 6

$ws.Range("C45").Value = 7
$ws.Range("C46").Value = 8
$ws.Range("C47").Value = 9

# --- D22/D23: change the text "100\n" to "7\n" while keeping the existing
#     wrap-text cell style (direct .Value assignment auto-coerces a
#     trim-numeric string like "7`n" into a plain number, losing both the
#     trailing newline and the text type, so build it via a scratch-cell
#     formula + paste-values round-trip, which keeps it text). ---
$ws.Range("Z100").Formula = "=""7""&CHAR(10)"
$ws.Range("Z100").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("Z100").Clear()

# --- Remove the now-obsolete 5th test block (rows 51-56) ---
$ws.Range("B51:D56").EntireRow.Delete()

# --- New column C width ---
$ws.Columns.Item(3).ColumnWidth = 29.833333333333332

# --- Selection moves to H22 (also drops the stale topLeftCell scroll anchor) ---
$ws.Range("H22").Select()
